# Remove the trailing empty paragraph at the end of the document
# (the empty paragraph that sits right before the section properties).
$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$prev = $d.Paragraphs.Item($count - 1)

# Delete from the end of the previous paragraph's mark through the end of
# the last (empty) paragraph's mark. This removes the last paragraph's
# mark entirely, merging it away so the previous paragraph (ending in ".")
# becomes the final paragraph in the document body.
$r = $d.Range($prev.Range.End - 1, $last.Range.End)
$r.Delete()
